$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '25.342.67'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -2.08%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.660.65'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -4.09%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9984'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.23%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '235.92'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -3.82%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9988'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.22%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4794'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -4.44%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2597'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -4.57%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06146'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -0.28%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07071'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -2.27%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.654.06'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -5.00%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '14.65'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -3.17%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.5848'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -10.25%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.359'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -8.60%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '74.14'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -3.67%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.9979'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -0.07%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.9993'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -0.18%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '25.369.16'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -2.04%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000006682'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -1.97%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.39'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -4.37%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.864.73'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -5.04%  '

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -4.76%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.584'

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.308'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -3.08%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '134.18'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +0.20%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '15.07'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -1.37%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.380'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -2.99%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '104.73'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -0.84%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.671'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -6.65%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '3.991'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +0.36%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.07642'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -6.00%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.598'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -2.34%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.04360'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -7.75%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.9976'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -0.07%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.601'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -1.91%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.6011'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -1.90%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.9395'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -5.58%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.624'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -4.48%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.8533'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -3.09%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.9990'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -0.10%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.01501'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -6.83%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '98.60'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -2.77%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.809'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -7.41%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.3732'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -4.30%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.669'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -6.79%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.1109'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -5.81%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '6.174'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -2.84%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.05242'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -0.70%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '29.36'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -4.57%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.213'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -1.63%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.000'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -0.02%  '
